$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values: reduce stored precision (custom accuracy rounding) ---
$ws.Range("B5").Value  = 7.23
$ws.Range("C5").Value  = 5.8
$ws.Range("D5").Value  = 0.16
$ws.Range("E5").Value  = 16.3
$ws.Range("F5").Value  = 13.49
$ws.Range("G5").Value  = 6.04
$ws.Range("H5").Value  = 21.85
$ws.Range("I5").Value  = 9.390000000000001
$ws.Range("J5").Value  = 4.28
$ws.Range("K5").Value  = 6.58
$ws.Range("L5").Value  = 7.06
$ws.Range("M5").Value  = 6.93
$ws.Range("N5").Value  = 1.81
$ws.Range("O5").Value  = 5.76
$ws.Range("P5").Value  = 8.529999999999999
$ws.Range("Q5").Value  = 4.96
$ws.Range("R5").Value  = 0.19
$ws.Range("S5").Value  = 0.15
$ws.Range("T5").Value  = 84.05
$ws.Range("U5").Value  = 16.76
$ws.Range("V5").Value  = 5.43
$ws.Range("W5").Value  = 11.02
$ws.Range("X5").Value  = 6.15
$ws.Range("Y5").Value  = 0.76
$ws.Range("Z5").Value  = 11.35
$ws.Range("AA5").Value = 4.9
$ws.Range("AB5").Value = 4.76
$ws.Range("AC5").Value = 5.38
$ws.Range("AD5").Value = 7.38
$ws.Range("AE5").Value = 0.13
$ws.Range("AF5").Value = 19.81
$ws.Range("AG5").Value = 3.33
$ws.Range("AH5").Value = 6.78

# --- Remove the last data row (row 6) entirely; shrinks dimension to A1:AH5 ---
$ws.Rows.Item(6).Delete()

# --- Narrow columns W (23) and Z (26) from width 8 to width 7 ---
$ws.Columns.Item(23).ColumnWidth = 6.17
$ws.Columns.Item(26).ColumnWidth = 6.17
